$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B width (Weekly Progress column width adjustment) ---
$ws.Columns.Item(2).ColumnWidth = 21

# --- Add row 8 (2025-06-02) ---
$a8 = $ws.Cells.Item(8, 1)
$a8.Value2 = 45810
$a8.NumberFormat = $ws.Cells.Item(7, 1).NumberFormat

$ws.Cells.Item(8, 2).Value2 = "Aanwezig - Documentation"
$ws.Cells.Item(8, 3).Value2 = "Aanwezig "
$ws.Cells.Item(8, 4).Value2 = "Aanwezig"
$ws.Cells.Item(8, 5).Value2 = "Afwezig"
$ws.Cells.Item(8, 6).Value2 = "Aanwezig - gepraat over groeps process"
$ws.Cells.Item(8, 7).Value2 = "Aanwezig"

# --- Add row 9 (2025-06-03) ---
$a9 = $ws.Cells.Item(9, 1)
$a9.Value2 = 45811
$a9.NumberFormat = $ws.Cells.Item(7, 1).NumberFormat

$ws.Cells.Item(9, 2).Value2 = "Aanwezig - Documentation"
$ws.Cells.Item(9, 3).Value2 = "Aanwezig - Demo"
$ws.Cells.Item(9, 4).Value2 = "Aanwezig"
$ws.Cells.Item(9, 5).Value2 = "Afwezig"
$ws.Cells.Item(9, 6).Value2 = "Aanwezig - groeps process vertellen"
$ws.Cells.Item(9, 7).Value2 = "Aanwezig"

# --- Update the view/selection state ---
[void]$ws.Range("G14").Select()

Write-Host "Edit applied successfully"
